$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for the new row 53 (date column uses the same
# custom date style as the other rows) by copying row 52's format down.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)

# Full replacement data for rows 2-53 (date_of_forecast, y_0, y_0_forecast,
# y_1, y_1_forecast). A new observation was inserted at the top (date
# 39400) shifting all subsequent rows down by one, a new row was appended
# at the bottom (date 45891), and the forecast values were recomputed.
$data = @(
    @(39400, 2007, 0.4235526809466261, 2008, 1.652713177926435),
    @(39583, 2008, -0.1474887606676356, 2009, 1.265019766896436),
    @(39765, 2008, -0.5718076928962645, 2009, 0.4618648366506939),
    @(39948, 2009, 0.337821977117625, 2010, 1.858657482882586),
    @(40130, 2009, 0.3486139762225005, 2010, -2.340608900318997),
    @(40310, 2010, -1.890773121057054, 2011, -3.246097549514837),
    @(40494, 2010, -0.1384957661262898, 2011, 1.659950937631938),
    @(40676, 2011, 1.89156560755015, 2012, 1.609625625599986),
    @(40862, 2011, 1.566479473280147, 2012, 2.497560574935442),
    @(41044, 2012, 1.113165545862094, 2013, 1.609625625600009),
    @(41228, 2012, 0.7307568962936939, 2013, 1.029202372425875),
    @(41409, 2013, 1.070385798714391, 2014, 3.238605209599998),
    @(41592, 2013, 0.818818812164257, 2014, 0.7004752402133052),
    @(41774, 2014, 1.384186838979806, 2015, 2.777885851461503),
    @(41957, 2014, 0.9180054319587239, 2015, 2.577405783391451),
    @(42137, 2015, 2.349355943833098, 2016, 2.436566844071941),
    @(42321, 2015, 1.984684278296656, 2016, 2.047428048848832),
    @(42503, 2016, 1.786425635558397, 2017, 1.694971351092267),
    @(42689, 2016, 1.755995812646982, 2017, 1.55296524673576),
    @(42867, 2017, 1.331333081915509, 2018, 1.216098605743365),
    @(43053, 2017, 1.946965557828384, 2018, 0.2326081529569146),
    @(43145, 2018, 0.522170282006873, 2019, 0.03619160485377471),
    @(43235, 2018, 1.282262557986447, 2019, 1.784618024189011),
    @(43326, 2018, 1.238324979098038, 2019, 1.281608622679231),
    @(43418, 2018, 1.06432145354225, 2019, -0.6993904531539141),
    @(43510, 2019, 0.5141202228481356, 2020, 0.8660061896410332),
    @(43600, 2019, 2.247109253368285, 2020, 4.887093273600018),
    @(43691, 2019, 1.254519294658696, 2020, 0.5495555957892417),
    @(43783, 2019, 1.361817904277696, 2020, 1.656460003703519),
    @(43875, 2020, 0.631197969589059, 2021, 0.07482640125562323),
    @(43966, 2020, -4.247034401476801, 2021, -12.1986023424),
    @(44068, 2020, -5.210209911466245, 2021, 17.92194299430632),
    @(44159, 2020, -4.352425014431304, 2021, -0.9164352768978534),
    @(44251, 2021, -5.694190915127639, 2022, -9.259896422400004),
    @(44341, 2021, -2.608215948579529, 2022, 12.1815000816919),
    @(44432, 2021, -1.897775264882628, 2022, 22.87548668341248),
    @(44525, 2021, -1.761645650979182, 2022, 2.211960525313206),
    @(44617, 2022, 3.750424156851184, 2023, 1.694448575331275),
    @(44706, 2022, 4.863085601670813, 2023, 5.870037016039187),
    @(44798, 2022, 4.787836378515364, 2023, 0.1238442324989597),
    @(44890, 2022, 5.20787683103745, 2023, 0.1996286546173343),
    @(44981, 2023, -1.15781484339289, 2024, -2.927305789734247),
    @(45071, 2023, -1.44371442952016, 2024, -1.24582517146522),
    @(45163, 2023, -0.49899188013105, 2024, 0.5969365658798464),
    @(45254, 2023, -0.9008525709169546, 2024, -0.01592365107300209),
    @(45345, 2024, 0.3201030582971054, 2025, 0.4006004000999486),
    @(45436, 2024, 0.06625622369935691, 2025, -0.1040473946152809),
    @(45534, 2024, 0.5084754301873051, 2025, 0.2859187804512109),
    @(45618, 2024, 0.2738544794132824, 2025, 0.7861875694807674),
    @(45713, 2025, 0.3933715771899493, 2026, 0.4005912323909122),
    @(45800, 2025, 0.9919038146506631, 2026, 0.8749339604052775),
    @(45891, 2025, 0.9584581489103794, 2026, 0.3060804787961446)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
